$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.958.90"
$ws.Range("E2").Value = "  +3.16%  "

# Row 3
$ws.Range("D3").Value = "1.777.09"
$ws.Range("E3").Value = "  -0.21%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9976"
$ws.Range("E4").Value = "  -0.78%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.32"
$ws.Range("E5").Value = "  -0.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9945"
$ws.Range("E6").Value = "  -0.74%  "

# Row 7
$ws.Range("E7").Value = "  -0.20%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3418"
$ws.Range("E8").Value = "  +0.42%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.89"
$ws.Range("E9").Value = "  -0.43%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.143"
$ws.Range("E10").Value = "  -3.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07430"
$ws.Range("E11").Value = "  +0.12%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.86"
$ws.Range("E12").Value = "  +5.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9954"
$ws.Range("E13").Value = "  -0.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.380"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15
$ws.Range("D15").Value = "1.774.54"
$ws.Range("E15").Value = "  -0.26%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.122"
$ws.Range("E16").Value = "  +1.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001079"
$ws.Range("E17").Value = "  -0.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06665"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.61"
$ws.Range("E19").Value = "  -0.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9962"
$ws.Range("E20").Value = "  -0.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.44"
$ws.Range("E21").Value = "  +1.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.427"
$ws.Range("E22").Value = "  -1.62%  "

# Row 23
$ws.Range("D23").Value = "27.938.10"
$ws.Range("E23").Value = "  +3.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.10"
$ws.Range("E24").Value = "  -0.76%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.372"
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.438"
$ws.Range("E26").Value = "  -1.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.78"
$ws.Range("E27").Value = "  -1.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.436"
$ws.Range("E28").Value = "  -2.55%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.00"
$ws.Range("E29").Value = "  -0.86%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.56"
$ws.Range("E30").Value = "  +0.47%  "

# Row 31
$ws.Range("B31").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C31").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D31").Value = "1.971.74"
$ws.Range("E31").Value = "  -0.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.184"
$ws.Range("E32").Value = "  +2.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.955"
$ws.Range("E33").Value = "  -0.75%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08794"
$ws.Range("E34").Value = "  +1.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.83"
$ws.Range("E35").Value = "  -1.64%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02437"
$ws.Range("E36").Value = "  +5.40%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6881"
$ws.Range("E37").Value = "  +0.94%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.341"
$ws.Range("E38").Value = "  -0.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06349"
$ws.Range("E39").Value = "  +1.14%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2193"
$ws.Range("E40").Value = "  +1.06%  "

# Row 41
$ws.Range("E41").Value = "  -6.56%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.238"
$ws.Range("E42").Value = "  +0.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.282"
$ws.Range("E43").Value = "  -2.95%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.26"
$ws.Range("E44").Value = "  +0.70%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9950"
$ws.Range("E45").Value = "  -0.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6310"
$ws.Range("E46").Value = "  -1.54%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.833"
$ws.Range("E47").Value = "  -0.58%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.21"
$ws.Range("E48").Value = "  +0.89%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.089"
$ws.Range("E49").Value = "  -1.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07374"
$ws.Range("E50").Value = "  +4.09%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.63"
$ws.Range("E51").Value = "  +0.23%  "
